# Updated cryptos list on Wed Dec 20 14:13:11 UTC 2023 with GitHub Actions
#
# Refreshes the Price/Volume(1h) columns of the crypto table with new
# quotes, and corrects the Litecoin/Uniswap (rows 20-21) and
# Cronos/WOONetwork (rows 47-48) row ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) cells whose new values would otherwise be
# auto-interpreted by Excel as numbers, so they remain text like the source data.
$textForceRows = @(5,6,9,10,11,12,17,20,21,22,23,24,26,27,31,32,33,34,36,37,39,40,41,42,43,44,45,46,47,48,51)
foreach ($r in $textForceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.750.59"
$ws.Range("E2").Value = "  +2.37%  "

$ws.Range("D3").Value = "2.241.41"
$ws.Range("E3").Value = "  +1.04%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "259.29"
$ws.Range("E5").Value = "  +3.30%  "

$ws.Range("D6").Value = "79.16"
$ws.Range("E6").Value = "  +5.99%  "

$ws.Range("E7").Value = "  +1.60%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  +2.53%  "

$ws.Range("D10").Value = "43.30"
$ws.Range("E10").Value = "  +6.68%  "

$ws.Range("D11").Value = "0.0927"
$ws.Range("E11").Value = "  +1.42%  "

$ws.Range("D12").Value = "7.10"
$ws.Range("E12").Value = "  +4.14%  "

$ws.Range("E13").Value = "  +1.75%  "

$ws.Range("D14").Value = "2.583.67"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("E15").Value = "  +2.29%  "

$ws.Range("D16").Value = "2.226.30"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").Value = "0.796"
$ws.Range("E17").Value = "  +1.67%  "

$ws.Range("D18").Value = "43.675.08"
$ws.Range("E18").Value = "  +2.38%  "

$ws.Range("E19").Value = "  +2.06%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "6.08"
$ws.Range("E20").Value = "  +3.07%  "

$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "71.56"
$ws.Range("E21").Value = "  +0.74%  "

$ws.Range("D22").Value = "2.32"
$ws.Range("E22").Value = "  +7.19%  "

$ws.Range("D23").Value = "233.55"
$ws.Range("E23").Value = "  +1.88%  "

$ws.Range("D24").Value = "9.43"
$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").Value = "42.24"
$ws.Range("E26").Value = "  +9.41%  "

$ws.Range("D27").Value = "10.85"
$ws.Range("E27").Value = "  +1.06%  "

$ws.Range("E28").Value = "  -1.61%  "

$ws.Range("E29").Value = "  +1.12%  "

$ws.Range("E30").Value = "  -1.07%  "

$ws.Range("D31").Value = "173.12"
$ws.Range("E31").Value = "  +1.97%  "

$ws.Range("D32").Value = "20.60"
$ws.Range("E32").Value = "  +2.52%  "

$ws.Range("D33").Value = "0.0877"
$ws.Range("E33").Value = "  +11.20%  "

$ws.Range("D34").Value = "5.30"
$ws.Range("E34").Value = "  +2.38%  "

$ws.Range("E35").Value = "  +1.29%  "

$ws.Range("D36").Value = "0.0366"
$ws.Range("E36").Value = "  +13.92%  "

$ws.Range("D37").Value = "4.49"
$ws.Range("E37").Value = "  +2.61%  "

$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").Value = "13.21"
$ws.Range("E39").Value = "  +10.41%  "

$ws.Range("D40").Value = "2.88"
$ws.Range("E40").Value = "  +20.14%  "

$ws.Range("D41").Value = "2.14"
$ws.Range("E41").Value = "  +2.73%  "

$ws.Range("D42").Value = "0.204"
$ws.Range("E42").Value = "  +1.32%  "

$ws.Range("D43").Value = "61.88"
$ws.Range("E43").Value = "  +5.10%  "

$ws.Range("D44").Value = "5.39"
$ws.Range("E44").Value = "  +2.05%  "

$ws.Range("D45").Value = "104.54"
$ws.Range("E45").Value = "  +2.14%  "

$ws.Range("D46").Value = "8.55"
$ws.Range("E46").Value = "  +0.30%  "

$ws.Range("B47").Value = "WOONetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D47").Value = "0.471"
$ws.Range("E47").Value = "  +0.38%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.0986"
$ws.Range("E48").Value = "  +0.43%  "

$ws.Range("E49").Value = "  +1.81%  "

$ws.Range("E50").Value = "  +2.15%  "

$ws.Range("D51").Value = "1.50"
$ws.Range("E51").Value = "  +25.52%  "
